# Scheduled-runner refresh of computed market/profit columns (H:N) across
# several leve-profit worksheets. Values come from an updated market-price
# pull; row identity (A:G) is unchanged.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 9705.388999999999
$ws.Range("I32").Value = 10074.875
$ws.Range("K32").Value = 10074.875
$ws.Range("M32").Value = -9748.875

$ws.Range("H33").Value = 3846269.5
$ws.Range("I33").Value = 4166785.8
$ws.Range("K33").Value = 4166785.8
$ws.Range("M33").Value = -4166556.8

$ws.Range("H43").Value = 14869.533
$ws.Range("I43").Value = 26077.5
$ws.Range("J43").Value = 10793.909
$ws.Range("K43").Value = 26077.5
$ws.Range("L43").Value = 10793.909
$ws.Range("M43").Value = -26008.5
$ws.Range("N43").Value = -10931.909

$ws.Range("H112").Value = 1770.909
$ws.Range("I112").Value = 996.3333
$ws.Range("J112").Value = 2061.375
$ws.Range("K112").Value = 2988.9999
$ws.Range("L112").Value = 6184.125
$ws.Range("M112").Value = -1880.9999
$ws.Range("N112").Value = -8400.125

$ws.Range("H116").Value = 13172
$ws.Range("I116").Value = 12910
$ws.Range("K116").Value = 12910
$ws.Range("M116").Value = -9468

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("L107").ClearContents()
$ws.Range("N107").Value = 0

$ws.Range("H112").Value = 51629
$ws.Range("J112").Value = 51629
$ws.Range("L112").Value = 51629
$ws.Range("N112").Value = -54583

$ws.Range("H121").Value = 58333.332
$ws.Range("J121").Value = 58333.332
$ws.Range("L121").Value = 58333.332
$ws.Range("N121").Value = -61827.332

$ws.Range("H133").Value = 297999.34
$ws.Range("I133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("M133").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 27496.5
$ws.Range("J82").Value = 80000
$ws.Range("L82").Value = 80000
$ws.Range("N82").Value = -80766

$ws.Range("H85").Value = 27496.5
$ws.Range("J85").Value = 80000
$ws.Range("L85").Value = 80000
$ws.Range("N85").Value = -82652

$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").ClearContents()
$ws.Range("N112").Value = 0

$ws.Range("H132").Value = 72818.664
$ws.Range("J132").Value = 74321.125
$ws.Range("L132").Value = 74321.125
$ws.Range("N132").Value = -84441.125

$ws.Range("H133").Value = 65438
$ws.Range("J133").Value = 65438
$ws.Range("L133").Value = 65438
$ws.Range("N133").Value = -75558

$ws.Range("H134").Value = 1760.2
$ws.Range("I134").Value = 1391.8182
$ws.Range("K134").Value = 4175.4546
$ws.Range("M134").Value = -1640.4546

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 80645.71000000001
$ws.Range("I31").Value = 86145.266
$ws.Range("J31").Value = 28400
$ws.Range("K31").Value = 86145.266
$ws.Range("L31").Value = 28400
$ws.Range("M31").Value = -85850.266
$ws.Range("N31").Value = -28990

$ws.Range("H34").Value = 80645.71000000001
$ws.Range("I34").Value = 86145.266
$ws.Range("J34").Value = 28400
$ws.Range("K34").Value = 86145.266
$ws.Range("L34").Value = 28400
$ws.Range("M34").Value = -85943.266
$ws.Range("N34").Value = -28804

$ws.Range("H54").Value = 34307.69
$ws.Range("I54").Value = 28000
$ws.Range("J54").Value = 34833.332
$ws.Range("K54").Value = 28000
$ws.Range("L54").Value = 34833.332
$ws.Range("M54").Value = -27342
$ws.Range("N54").Value = -36149.332

$ws.Range("H58").Value = 1894.4286
$ws.Range("I58").Value = 1725.421
$ws.Range("K58").Value = 1725.421
$ws.Range("M58").Value = -1522.421

$ws.Range("H62").Value = 5469.4614
$ws.Range("I62").Value = 5645.727
$ws.Range("K62").Value = 5645.727
$ws.Range("M62").Value = -5021.727

$ws.Range("H65").Value = 5469.4614
$ws.Range("I65").Value = 5645.727
$ws.Range("K65").Value = 28228.635
$ws.Range("M65").Value = -25108.635

$ws.Range("H132").Value = 3221.9678
$ws.Range("I132").Value = 3163.6072
$ws.Range("J132").Value = 3766.6667
$ws.Range("K132").Value = 9490.821599999999
$ws.Range("L132").Value = 11300.0001
$ws.Range("M132").Value = -6960.821599999999
$ws.Range("N132").Value = -16360.0001

$ws.Range("H133").Value = 47800
$ws.Range("J133").Value = 47800
$ws.Range("L133").Value = 47800
$ws.Range("N133").Value = -52860

$ws.Range("H136").Value = 1894.4286
$ws.Range("I136").Value = 1725.421
$ws.Range("K136").Value = 5176.263
$ws.Range("M136").Value = -2626.263

$ws.Range("H137").Value = 49714.285
$ws.Range("J137").Value = 49714.285
$ws.Range("L137").Value = 49714.285
$ws.Range("N137").Value = -59914.285

$ws.Range("H139").Value = 142000
$ws.Range("J139").Value = 142000
$ws.Range("L139").Value = 142000
$ws.Range("N139").Value = -152280

$ws.Range("H140").Value = 100707
$ws.Range("J140").Value = 100707
$ws.Range("L140").Value = 100707
$ws.Range("N140").Value = -111067

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 3000
$ws.Range("J33").Value = 3000
$ws.Range("L33").Value = 3000
$ws.Range("N33").Value = -3504

$ws.Range("H70").Value = 10029.538
$ws.Range("I70").Value = 11923
$ws.Range("K70").Value = 11923
$ws.Range("M70").Value = -11653

$ws.Range("H73").Value = 10029.538
$ws.Range("I73").Value = 11923
$ws.Range("K73").Value = 11923
$ws.Range("M73").Value = -10987

$ws.Range("H102").Value = 2403.9443
$ws.Range("I102").Value = 1663.7059
$ws.Range("J102").Value = 14988
$ws.Range("K102").Value = 1663.7059
$ws.Range("L102").Value = 14988
$ws.Range("M102").Value = -41.70589999999993
$ws.Range("N102").Value = -18232

$ws.Range("H107").Value = 1452
$ws.Range("I107").Value = 1270.4286
$ws.Range("K107").Value = 1270.4286
$ws.Range("M107").Value = 649.5714

$ws.Range("H132").Value = 219813.72
$ws.Range("I132").Value = 224587.36
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 673762.08
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -671232.08
$ws.Range("N132").Value = -20060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2840.9211
$ws.Range("I61").Value = 2567.0334
$ws.Range("J61").Value = 3868
$ws.Range("K61").Value = 2567.0334
$ws.Range("L61").Value = 3868
$ws.Range("M61").Value = -2365.0334
$ws.Range("N61").Value = -4272

$ws.Range("H68").Value = 5397.95
$ws.Range("I68").Value = 3886.4443
$ws.Range("J68").Value = 19001.5
$ws.Range("K68").Value = 3886.4443
$ws.Range("L68").Value = 19001.5
$ws.Range("M68").Value = -3137.4443
$ws.Range("N68").Value = -20499.5

$ws.Range("H71").Value = 5397.95
$ws.Range("I71").Value = 3886.4443
$ws.Range("J71").Value = 19001.5
$ws.Range("K71").Value = 19432.2215
$ws.Range("L71").Value = 95007.5
$ws.Range("M71").Value = -15688.2215
$ws.Range("N71").Value = -102495.5

$ws.Range("H100").Value = 3037.5557
$ws.Range("I100").Value = 3037.5557
$ws.Range("K100").Value = 3037.5557
$ws.Range("M100").Value = -2496.5557

$ws.Range("H110").Value = 28624.75
$ws.Range("J110").Value = 28624.75
$ws.Range("L110").Value = 28624.75
$ws.Range("N110").Value = -36804.75

$ws.Range("H113").Value = 2840.9211
$ws.Range("I113").Value = 2567.0334
$ws.Range("J113").Value = 3868
$ws.Range("K113").Value = 2567.0334
$ws.Range("L113").Value = 3868
$ws.Range("M113").Value = -397.0333999999998
$ws.Range("N113").Value = -8208

$ws.Range("H127").Value = 66432.5
$ws.Range("J127").Value = 66432.5
$ws.Range("L127").Value = 66432.5
$ws.Range("N127").Value = -76352.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 19096.04
$ws.Range("I62").Value = 18590.863
$ws.Range("J62").Value = 22271.428
$ws.Range("K62").Value = 18590.863
$ws.Range("L62").Value = 22271.428
$ws.Range("M62").Value = -17966.863
$ws.Range("N62").Value = -23519.428

$ws.Range("H65").Value = 19096.04
$ws.Range("I65").Value = 18590.863
$ws.Range("J65").Value = 22271.428
$ws.Range("K65").Value = 92954.315
$ws.Range("L65").Value = 111357.14
$ws.Range("M65").Value = -89834.315
$ws.Range("N65").Value = -117597.14
